$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 updates
$ws.Range("G2").Value = 1.38
$ws.Range("H2").Value = 4.33
$ws.Range("I2").Value = 9.5
$ws.Range("J2").Value = 1.95
$ws.Range("L2").Value = 9
$ws.Range("M2").Value = 1.08
$ws.Range("N2").Value = 8
$ws.Range("S2").Value = 2.1
$ws.Range("T2").Value = 1.7
$ws.Range("W2").Value = 3.75
$ws.Range("X2").Value = 1.25
$ws.Range("AE2").Value = 9.5
$ws.Range("AK2").Value = 29
$ws.Range("AO2").Value = 29
$ws.Range("AP2").Value = 126
$ws.Range("AQ2").Value = 81

# Row 3 updates
$ws.Range("G3").Value = 3.4
$ws.Range("H3").Value = 3
$ws.Range("M3").Value = 1.11
$ws.Range("N3").Value = 6.5
$ws.Range("U3").Value = 3.9
$ws.Range("Y3").Value = 1.53
$ws.Range("Z3").Value = 2.38
$ws.Range("AA3").Value = 2
$ws.Range("AB3").Value = 1.73
$ws.Range("AC3").Value = 8.5
$ws.Range("AF3").Value = 41
$ws.Range("AI3").Value = 6.5
$ws.Range("AK3").Value = 17
$ws.Range("AM3").Value = 6.5
$ws.Range("AN3").Value = 10

# Row 4 updates
$ws.Range("G4").Value = 2.55
$ws.Range("I4").Value = 2.75
$ws.Range("J4").Value = 3.4
$ws.Range("L4").Value = 3.6
$ws.Range("O4").Value = 1.44
$ws.Range("P4").Value = 2.63
$ws.Range("Q4").Value = 1.8
$ws.Range("R4").Value = 2
$ws.Range("S4").Value = 2.4
$ws.Range("T4").Value = 1.53
$ws.Range("W4").Value = 4.5
$ws.Range("X4").Value = 1.18
$ws.Range("AC4").Value = 7
$ws.Range("AF4").Value = 26
$ws.Range("AG4").Value = 23
$ws.Range("AL4").Value = 67
$ws.Range("AM4").Value = 7.5
$ws.Range("AN4").Value = 13
$ws.Range("AP4").Value = 29
$ws.Range("AQ4").Value = 26
